# Disponibilidad.xlsx - "Update automatico via Actualizar 02-05-2021 01-12-57"
#
# 1) The 14 rows of the previous check-run (rows 240-253) get their
#    timestamp (column D) corrected from 44232.02931486353 to
#    44232.02931486111 (same moment, re-serialised float).
# 2) A brand-new check-run of 14 rows (254-267) is appended, duplicating
#    the same Name/URL/Disponible pattern with a fresh timestamp
#    (44232.05062243481) and fresh external hyperlinks on column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix the stored timestamp for the existing last block (rows 240-253)
# ---------------------------------------------------------------------
for ($r = 240; $r -le 253; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.02931486111
}

# ---------------------------------------------------------------------
# 2) Append the new block (rows 254-267) as a duplicate of 240-253,
#    then touch-up the timestamp and hyperlinks for the new rows.
# ---------------------------------------------------------------------
$ws.Range("A240:D253").Copy($ws.Range("A254:D267"))

$targets = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

for ($i = 0; $i -lt 14; $i++) {
    $r = 254 + $i
    $ws.Cells.Item($r, 4).Value = 44232.05062243481

    if ($r -eq 262) {
        # MapStore row keeps the trailing "#/" as a hyperlink sub-address
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $targets[$i], "/") | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $targets[$i]) | Out-Null
    }

    # Hyperlinks.Add() re-stamps the anchor cell with a brand-new style;
    # put it back on the shared "Hyperlink" cell style used by every
    # other URL cell in column B.
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
}
